$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "2023-11-06T00:00:00"
$ws.Range("C5").Select()
